$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1995.9524
$ws.Range("I28").Value = 2080.75
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 2080.75
$ws.Range("L28").Value = 300
$ws.Range("M28").Value = -1595.75
$ws.Range("N28").Value = -1270
$ws.Range("H40").Value = 1980
$ws.Range("I40").Value = 1881.8182
$ws.Range("K40").Value = 1881.8182
$ws.Range("M40").Value = -1706.8182
$ws.Range("H98").Value = 1893.5344
$ws.Range("I98").Value = 2018.0577
$ws.Range("J98").Value = 814.3333
$ws.Range("K98").Value = 2018.0577
$ws.Range("L98").Value = 814.3333
$ws.Range("M98").Value = -520.0577000000001
$ws.Range("N98").Value = -3810.3333
$ws.Range("H101").Value = 727.7143
$ws.Range("I101").Value = 727.7143
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2183.1429
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -561.1428999999998
$ws.Range("N101").Value = $null
$ws.Range("H116").Value = 3065.9167
$ws.Range("I116").Value = 2177
$ws.Range("J116").Value = 3700.8572
$ws.Range("K116").Value = 2177
$ws.Range("L116").Value = 3700.8572
$ws.Range("M116").Value = 1265
$ws.Range("N116").Value = -10584.8572
$ws.Range("H122").Value = 1893.5344
$ws.Range("I122").Value = 2018.0577
$ws.Range("J122").Value = 814.3333
$ws.Range("K122").Value = 6054.1731
$ws.Range("L122").Value = 2442.9999
$ws.Range("M122").Value = -3604.1731
$ws.Range("N122").Value = -7342.9999
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H132").Value = 12354623
$ws.Range("I132").Value = 14499788
$ws.Range("J132").Value = 19926.25
$ws.Range("K132").Value = 43499364
$ws.Range("L132").Value = 59778.75
$ws.Range("M132").Value = -43496834
$ws.Range("N132").Value = -64838.75
$ws.Range("H137").Value = 1413.1555
$ws.Range("I137").Value = 1010.5
$ws.Range("J137").Value = 1798.3043
$ws.Range("K137").Value = 3031.5
$ws.Range("L137").Value = 5394.9129
$ws.Range("M137").Value = -481.5
$ws.Range("N137").Value = -10494.9129

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4636.4287
$ws.Range("I31").Value = 4636.4287
$ws.Range("K31").Value = 4636.4287
$ws.Range("M31").Value = -4342.4287
$ws.Range("H32").Value = 12838.2
$ws.Range("I32").Value = 11213.292
$ws.Range("K32").Value = 11213.292
$ws.Range("M32").Value = -10926.292
$ws.Range("H45").Value = 1205
$ws.Range("I45").Value = 1069.375
$ws.Range("J45").Value = 1566.6666
$ws.Range("K45").Value = 1069.375
$ws.Range("L45").Value = 1566.6666
$ws.Range("M45").Value = -692.375
$ws.Range("N45").Value = -2320.6666
$ws.Range("H61").Value = 142859400
$ws.Range("I61").Value = 333334600
$ws.Range("K61").Value = 333334600
$ws.Range("M61").Value = -333334388
$ws.Range("H136").Value = 142859400
$ws.Range("I136").Value = 333334600
$ws.Range("K136").Value = 1000003800
$ws.Range("M136").Value = -1000001250

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2815.2856
$ws.Range("I20").Value = 2288.4
$ws.Range("J20").Value = 4132.5
$ws.Range("K20").Value = 2288.4
$ws.Range("L20").Value = 4132.5
$ws.Range("M20").Value = -2041.4
$ws.Range("N20").Value = -4626.5
$ws.Range("H134").Value = 9116.866
$ws.Range("I134").Value = 1896.3
$ws.Range("J134").Value = 23558
$ws.Range("K134").Value = 5688.9
$ws.Range("L134").Value = 70674
$ws.Range("M134").Value = -3153.9
$ws.Range("N134").Value = -75744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1252.1
$ws.Range("I31").Value = 1064.7297
$ws.Range("J31").Value = 3563
$ws.Range("K31").Value = 1064.7297
$ws.Range("L31").Value = 3563
$ws.Range("M31").Value = -769.7297000000001
$ws.Range("N31").Value = -4153
$ws.Range("H34").Value = 1252.1
$ws.Range("I34").Value = 1064.7297
$ws.Range("J34").Value = 3563
$ws.Range("K34").Value = 1064.7297
$ws.Range("L34").Value = 3563
$ws.Range("M34").Value = -862.7297000000001
$ws.Range("N34").Value = -3967
$ws.Range("H94").Value = 1265.3077
$ws.Range("I94").Value = 1071
$ws.Range("J94").Value = 1431.8572
$ws.Range("K94").Value = 1071
$ws.Range("L94").Value = 1431.8572
$ws.Range("M94").Value = -620
$ws.Range("N94").Value = -2333.8572
$ws.Range("H122").Value = 893.6667
$ws.Range("I122").Value = 720.6667
$ws.Range("J122").Value = 1066.6666
$ws.Range("K122").Value = 2162.0001
$ws.Range("L122").Value = 3199.9998
$ws.Range("M122").Value = 287.9998999999998
$ws.Range("N122").Value = -8099.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 204.2381
$ws.Range("I14").Value = 204.2381
$ws.Range("K14").Value = 612.7143
$ws.Range("M14").Value = -439.7143
$ws.Range("H122").Value = 1081.0476
$ws.Range("J122").Value = 1157.125
$ws.Range("L122").Value = 10414.125
$ws.Range("N122").Value = -15314.125
$ws.Range("H137").Value = 27781040
$ws.Range("I137").Value = 125000664
$ws.Range("J137").Value = 4004.0952
$ws.Range("K137").Value = 375001992
$ws.Range("L137").Value = 12012.2856
$ws.Range("M137").Value = -374996892
$ws.Range("N137").Value = -22212.2856
$ws.Range("H140").Value = 31942.082
$ws.Range("I140").Value = 66307.586
$ws.Range("J140").Value = 2731.4
$ws.Range("K140").Value = 198922.758
$ws.Range("L140").Value = 8194.200000000001
$ws.Range("M140").Value = -193742.758
$ws.Range("N140").Value = -18554.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H82").Value = 2191.0715
$ws.Range("I82").Value = 2266.5557
$ws.Range("J82").Value = 2055.2
$ws.Range("K82").Value = 2266.5557
$ws.Range("L82").Value = 2055.2
$ws.Range("M82").Value = -1905.5557
$ws.Range("N82").Value = -2777.2
$ws.Range("H85").Value = 2191.0715
$ws.Range("I85").Value = 2266.5557
$ws.Range("J85").Value = 2055.2
$ws.Range("K85").Value = 2266.5557
$ws.Range("L85").Value = 2055.2
$ws.Range("M85").Value = -1018.5557
$ws.Range("N85").Value = -4551.2
$ws.Range("H122").Value = 14716311
$ws.Range("I122").Value = 20843808
$ws.Range("K122").Value = 62531424
$ws.Range("M122").Value = -62528974

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H132").Value = 2555.4546
$ws.Range("I132").Value = 1322.2
$ws.Range("J132").Value = 3583.1667
$ws.Range("K132").Value = 3966.6
$ws.Range("L132").Value = 10749.5001
$ws.Range("M132").Value = -1436.6
$ws.Range("N132").Value = -15809.5001
$ws.Range("H136").Value = 1036.2667
$ws.Range("I136").Value = 1037.2632
$ws.Range("J136").Value = 1034.5454
$ws.Range("K136").Value = 3111.7896
$ws.Range("L136").Value = 3103.6362
$ws.Range("M136").Value = -561.7896000000001
$ws.Range("N136").Value = -8203.636200000001
